$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (Брюн Феликс): fill C6:H6 with "ОК"/"ок" and add I6 = "ок" ---
$ws.Range("C6").Value = "ОК"
$ws.Range("D6").Value = "ок"
$ws.Range("E6").Value = "ОК"
$ws.Range("F6").Value = "ОК"
$ws.Range("G6").Value = "ОК"
$ws.Range("H6").Value = "ОК"

# I6 is a brand-new cell; copy the format from I9 (the other "ок"-styled cell)
# so it picks up the correct style (s="4") before writing its value.
$ws.Range("I9").Copy()
$ws.Range("I6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I6").Value = "ок"

# --- Row 11 (Егунов Захар): fill C11:I11 with "ОК" ---
$ws.Range("C11").Value = "ОК"
$ws.Range("D11").Value = "ОК"
$ws.Range("E11").Value = "ОК"
$ws.Range("F11").Value = "ОК"
$ws.Range("G11").Value = "ОК"
$ws.Range("H11").Value = "ОК"

# I11 is a brand-new cell; copy the format from I10 (a regular "ОК"-styled cell)
$ws.Range("I10").Copy()
$ws.Range("I11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I11").Value = "ОК"

# --- Row 19 (Никифорова Юлия): fill D19:I19 with "ОК" (C19 already "ОК") ---
$ws.Range("D19").Value = "ОК"
$ws.Range("E19").Value = "ОК"
$ws.Range("F19").Value = "ОК"
$ws.Range("G19").Value = "ОК"
$ws.Range("H19").Value = "ОК"

# I19 is a brand-new cell; copy the format from I18 (a regular "ОК"-styled cell)
$ws.Range("I18").Copy()
$ws.Range("I19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I19").Value = "ОК"

$excel.CutCopyMode = 0

# Matches the recorded selection change in the workbook after the edit.
$ws.Range("C6:I6").Select()
